# Q3 Update - 2025
# Applies the UNHCR "UN-AFG" dataset refresh for the year-2024 block:
#   - refresh the short-url token used throughout the sheet
#   - correct Afghanistan's 2024 idps/returned_idps/ooc/hst totals
#   - insert a new India row into the 2024 block (between Bangladesh and Iran)
#   - correct Iran's 2024 asylum_seekers total
#   - correct Pakistan's 2024 refugees/asylum_seekers totals
#   - renumber the sequential "items" counter for the rows that followed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Refresh the short-url value used by every data row (column B, rows 2-127)
$ws.Range("B2:B127").Value = "MGuy8C"

# 2. Afghanistan / 2024 row (row 123): update idps, returned_idps, ooc, hst
$ws.Range("Q123").Value = 3199710
$ws.Range("R123").Value = 22687
$ws.Range("T123").Value = 72886
$ws.Range("V123").Value = 600000

# 3. Insert a new row for India right after Bangladesh (row 124), shifting
#    Iran/Pakistan/Tajikistan down by one row.
$ws.Rows.Item(125).Insert()

$ws.Range("A125").Value = 1
$ws.Range("B125").Value = "MGuy8C"
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 124
$ws.Range("E125").Value = 2024
$ws.Range("F125").Value = 88
$ws.Range("G125").Value = "India"
$ws.Range("H125").Value = "IND"
$ws.Range("I125").Value = "IND"
$ws.Range("J125").Value = 2
$ws.Range("K125").Value = "Afghanistan"
$ws.Range("L125").Value = "AFG"
$ws.Range("M125").Value = "AFG"
$ws.Range("N125").Value = 0
$ws.Range("O125").Value = 5
$ws.Range("P125").Value = 0
$ws.Range("Q125").Value = 0
$ws.Range("R125").Value = 0
$ws.Range("S125").Value = 0
$ws.Range("T125").Value = 0
$ws.Range("U125").Value = "-"
$ws.Range("V125").Value = 0

# 4. Renumber the "items" sequence for the rows pushed down by the insert
#    (Iran, Pakistan, Tajikistan are now rows 126, 127, 128).
$ws.Range("D126").Value = 125
$ws.Range("D127").Value = 126
$ws.Range("D128").Value = 127

# 5. Iran (now row 126): asylum_seekers corrected
$ws.Range("O126").Value = 55

# 6. Pakistan (now row 127): refugees/asylum_seekers corrected
$ws.Range("N127").Value = 20827
$ws.Range("O127").Value = 299
